$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A updates
$ws.Range("A2").Value = "Sanity31"
$ws.Range("A3").Value = "Sanity32"

# Column B updates
$ws.Range("B2").Value = "Test31"
$ws.Range("B3").Value = "Test32"

# Column D updates (leading apostrophe preserves the existing quote-prefixed
# text style on these cells instead of Excel re-styling them on assignment)
$ws.Range("D2").Value = "'Sanity31.Test31@gmail.com"
$ws.Range("D3").Value = "'Sanity32.Test32@gmail.com"

# Columns E and G updates (same quote-prefix style preservation)
$ws.Range("E2").Value = "'30"
$ws.Range("G2").Value = "'30"

# Update the active selection/cell to match the final state
$ws.Range("G3").Select()
